# Update "想去人数" (F column) values for sheets "展览" (index 1) and
# "全部类型" (index 4). Both sheets hold mirrored event data; most rows
# get the same updated value on both sheets, but rows 11 and 36 diverge.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

# Column F = column index 6. Rows where both sheets end up with the same
# new value.
$commonUpdates = @{
    2  = 214
    3  = 5607
    5  = 685
    6  = 677
    7  = 33
    8  = 16
    12 = 5651
    14 = 295
    17 = 28
    18 = 118
    19 = 4593
    25 = 79
    26 = 209
    30 = 355
    31 = 41
    32 = 47
    34 = 21
    35 = 34
    37 = 47
}

foreach ($row in $commonUpdates.Keys) {
    $newVal = $commonUpdates[$row]
    $ws1.Cells.Item($row, 6).Value = $newVal
    $ws4.Cells.Item($row, 6).Value = $newVal
}

# Row 11 ends at the same value (1593) on both sheets even though the
# starting values differed (1588 vs 1589).
$ws1.Cells.Item(11, 6).Value = 1593
$ws4.Cells.Item(11, 6).Value = 1593

# Row 36 diverges between the two sheets.
$ws1.Cells.Item(36, 6).Value = 5
$ws4.Cells.Item(36, 6).Value = 6
